# Weekly price-sheet update for "Hortaliza, Vega Modelo de Temuco - Cilantro".
# Two new daily records are inserted at the top of the data block (rows 529-530),
# pushing the previously-existing rows 529:605 down to 531:607.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 529-530; everything from the old row 529 onward shifts
# down by two rows (old 529 -> 531, ..., old 605 -> 607). This also grows the
# sheet's used range from A1:R605 to A1:R607.
$ws.Range("529:530").EntireRow.Insert()

# New row 529
$ws.Range("A529").Value = 10
$ws.Range("B529").Value = 'Vega Modelo de Temuco'
$ws.Range("C529").Value = 'La Araucanía'
$ws.Range("D529").Value = 45127
$ws.Range("E529").Value = 9
$ws.Range("F529").Value = 100112040
$ws.Range("G529").Value = 'Cilantro'
$ws.Range("H529").Value = 'Sin especificar'
$ws.Range("I529").Value = 'Primera'
$ws.Range("J529").Value = 50
$ws.Range("K529").Value = 4000
$ws.Range("L529").Value = 4000
$ws.Range("M529").Value = 4000
$ws.Range("N529").Value = '$/docena de atados (2 kilos)'
$ws.Range("O529").Value = 'Provincia de Cautín'
$ws.Range("P529").Value = 2000
$ws.Range("Q529").Value = 2
$ws.Range("R529").Value = 'Hortaliza'

# New row 530
$ws.Range("A530").Value = 10
$ws.Range("B530").Value = 'Vega Modelo de Temuco'
$ws.Range("C530").Value = 'La Araucanía'
$ws.Range("D530").Value = 45127
$ws.Range("E530").Value = 9
$ws.Range("F530").Value = 100112040
$ws.Range("G530").Value = 'Cilantro'
$ws.Range("H530").Value = 'Sin especificar'
$ws.Range("I530").Value = 'Primera'
$ws.Range("J530").Value = 200
$ws.Range("K530").Value = 4000
$ws.Range("L530").Value = 4600
$ws.Range("M530").Value = 4360
$ws.Range("N530").Value = '$/docena de atados (2 kilos)'
$ws.Range("O530").Value = 'Región Metropolitana'
$ws.Range("P530").Value = 2180
$ws.Range("Q530").Value = 2
$ws.Range("R530").Value = 'Hortaliza'
